$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9919543333333333
$ws.Range("H2").Value = 2.975863
$ws.Range("I2").Value = 0.008811579445878926
$ws.Range("J2").Value = 0.008811579445878926
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.689526666666667
$ws.Range("N2").Value = 5.06858
$ws.Range("O2").Value = 0.3793973676518417
$ws.Range("P2").Value = 0.3793973676518417
$ws.Range("Q2").Value = 1.675933298282222
$ws.Range("R2").Value = 15.08339968454
$ws.Range("S2").Value = 0.003343090046621539
$ws.Range("T2").Value = 0.003343090046621538
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9919543333333333
$ws.Range("H3").Value = 2.975863
$ws.Range("I3").Value = 0.008811579445878926
$ws.Range("J3").Value = 0.008811579445878926
$ws.Range("O3").Value = 0.5332368480836376
$ws.Range("P3").Value = 0.5332368480836375
$ws.Range("Q3").Value = 2.355497074493444
$ws.Range("R3").Value = 21.199473670441
$ws.Range("S3").Value = 0.004698658850359044
$ws.Range("T3").Value = 0.004698658850359043
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9919543333333333
$ws.Range("H4").Value = 2.975863
$ws.Range("I4").Value = 0.008811579445878926
$ws.Range("J4").Value = 0.008811579445878926
$ws.Range("M4").Value = 0.389056
$ws.Range("N4").Value = 1.167168
$ws.Range("O4").Value = 0.0873657842645208
$ws.Range("P4").Value = 0.08736578426452078
$ws.Range("Q4").Value = 0.3859257851093333
$ws.Range("R4").Value = 3.473332065984
$ws.Range("S4").Value = 0.0007698305488983439
$ws.Range("T4").Value = 0.0007698305488983438
$ws.Range("I5").Value = 0.6711393126876655
$ws.Range("J5").Value = 0.6711393126876655
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.689526666666667
$ws.Range("N5").Value = 5.06858
$ws.Range("O5").Value = 0.3793973676518417
$ws.Range("P5").Value = 0.3793973676518417
$ws.Range("Q5").Value = 127.6484799153178
$ws.Range("R5").Value = 1148.83631923786
$ws.Range("S5").Value = 0.2546284885613666
$ws.Range("T5").Value = 0.2546284885613666
$ws.Range("I6").Value = 0.6711393126876655
$ws.Range("J6").Value = 0.6711393126876655
$ws.Range("O6").Value = 0.5332368480836376
$ws.Range("P6").Value = 0.5332368480836375
$ws.Range("S6").Value = 0.3578762117225896
$ws.Range("T6").Value = 0.3578762117225895
$ws.Range("I7").Value = 0.6711393126876655
$ws.Range("J7").Value = 0.6711393126876655
$ws.Range("M7").Value = 0.389056
$ws.Range("N7").Value = 1.167168
$ws.Range("O7").Value = 0.0873657842645208
$ws.Range("P7").Value = 0.08736578426452078
$ws.Range("Q7").Value = 29.39427236145066
$ws.Range("R7").Value = 264.548451253056
$ws.Range("S7").Value = 0.05863461240370935
$ws.Range("T7").Value = 0.05863461240370934
$ws.Range("G8").Value = 36.02919333333333
$ws.Range("H8").Value = 108.08758
$ws.Range("I8").Value = 0.3200491078664556
$ws.Range("J8").Value = 0.3200491078664556
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.689526666666667
$ws.Range("N8").Value = 5.06858
$ws.Range("O8").Value = 0.3793973676518417
$ws.Range("P8").Value = 0.3793973676518417
$ws.Range("Q8").Value = 60.87228291515555
$ws.Range("R8").Value = 547.8505462364
$ws.Range("S8").Value = 0.1214257890438536
$ws.Range("T8").Value = 0.1214257890438536
$ws.Range("G9").Value = 36.02919333333333
$ws.Range("H9").Value = 108.08758
$ws.Range("I9").Value = 0.3200491078664556
$ws.Range("J9").Value = 0.3200491078664556
$ws.Range("O9").Value = 0.5332368480836376
$ws.Range("P9").Value = 0.5332368480836375
$ws.Range("Q9").Value = 85.55500655745111
$ws.Range("R9").Value = 769.99505901706
$ws.Range("S9").Value = 0.1706619775106889
$ws.Range("T9").Value = 0.1706619775106889
$ws.Range("G10").Value = 36.02919333333333
$ws.Range("H10").Value = 108.08758
$ws.Range("I10").Value = 0.3200491078664556
$ws.Range("J10").Value = 0.3200491078664556
$ws.Range("M10").Value = 0.389056
$ws.Range("N10").Value = 1.167168
$ws.Range("O10").Value = 0.0873657842645208
$ws.Range("P10").Value = 0.08736578426452078
$ws.Range("Q10").Value = 14.01737384149333
$ws.Range("R10").Value = 126.15636457344
$ws.Range("S10").Value = 0.0279613413119131
$ws.Range("T10").Value = 0.0279613413119131
